$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 60.3858856367471
$ws.Range("K2").Value = 45.265798412196
$ws.Range("L2").Value = 105.041401045725
$ws.Range("N2").Value = 36.56805362158
$ws.Range("B3").Value = 4082.89927538263
$ws.Range("K3").Value = 10.3482752823576
$ws.Range("N3").Value = 10686.2995154378
$ws.Range("B4").Value = 4.43677621365884
$ws.Range("H4").Value = 4.93007252540713
$ws.Range("I4").Value = 8.20702036594838
$ws.Range("K4").Value = 2.11723414529968
$ws.Range("N4").Value = 4.20162921972641
$ws.Range("B5").Value = 2.63809340787633
$ws.Range("K5").Value = 0.950961050927102
$ws.Range("N5").Value = 1.47147866946322
$ws.Range("B6").Value = 10.24651219227
$ws.Range("D6").Value = 5.46899451517479
$ws.Range("K6").Value = 2.8142556576562
$ws.Range("N6").Value = 3.47685915933737
$ws.Range("B7").Value = 3.11260802520171
$ws.Range("K7").Value = 1.84740260667086
$ws.Range("N7").Value = 2.2643651802235
$ws.Range("B8").Value = 0.433242778030893
$ws.Range("K8").Value = 0.378545006272694
$ws.Range("N8").Value = 0.435232581680904
$ws.Range("B9").Value = 2.21887505415695
$ws.Range("K9").Value = 2.02663842919111
$ws.Range("N9").Value = 2.51963186388005
$ws.Range("B10").Value = 3.22395138306479
$ws.Range("K10").Value = 4.53654976738986
$ws.Range("N10").Value = 1.04044221625182
$ws.Range("B11").Value = 3.58967842130771
$ws.Range("K11").Value = 2.91027449070564
$ws.Range("N11").Value = 3.85123225935303
$ws.Range("B12").Value = 2995.49159409664
$ws.Range("K12").Value = 1.21021581454767
$ws.Range("N12").Value = 1.48587141864868
$ws.Range("B13").Value = 1.19296152731599
$ws.Range("K13").Value = 1.03122036027141
$ws.Range("M13").Value = 0.465609222021071
$ws.Range("N13").Value = 0.989872452795338
$ws.Range("B14").Value = 3.01610246617166
$ws.Range("H14").Value = 9.07972714272027
$ws.Range("K14").Value = 2.63804254248156
$ws.Range("N14").Value = 2.06264842705817
$ws.Range("B15").Value = 2996.48664996978
$ws.Range("K15").Value = 2.79948849986231
$ws.Range("N15").Value = 2.21337105152458
$ws.Range("B16").Value = 5.02715067817026
$ws.Range("K16").Value = 3.06876640759826
$ws.Range("N16").Value = 2.94328605247035
$ws.Range("B17").Value = 3.57891277014804
$ws.Range("K17").Value = 4.52648996429164
$ws.Range("N17").Value = 3.39553892417094
$ws.Range("B18").Value = 2995.46663046482
$ws.Range("K18").Value = 1.31842237561983
$ws.Range("N18").Value = 1.59889206546569
$ws.Range("B19").Value = -11217.5158065268
$ws.Range("C19").Value = -20690.7934752408
$ws.Range("D19").Value = 0.746733701835473
$ws.Range("K19").Value = 0.743016011052597
$ws.Range("L19").Value = 105.041401045725
$ws.Range("N19").Value = -10683.6818133583

Write-Output "applied 63 changes"